$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph. It is the paragraph
#    that starts with a bold "Meta description" run, right after the
#    "Play Gate of The Pharaohs Slot for Free" Heading1 paragraph.
# ---------------------------------------------------------------------
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Gate of The Pharaohs Slot for
#    Free" right before the closing italic "feature image" paragraph,
#    and swap that paragraph's long image-prompt text for the short
#    meta-description text (keeping its italic formatting).
# ---------------------------------------------------------------------
$oldPrompt = "Create a cartoon-style feature image for Gate of The Pharaohs with a happy Maya warrior wearing glasses. The image should feature the Maya warrior standing in front of a large, golden gate adorned with ancient Egyptian symbols, such as the Eye of Ra and hieroglyphics. The gate should have an opening that reveals a glimpse of the riches waiting inside. The background should be a desert landscape with pyramids in the distance. The Maya warrior should be holding a treasure chest overflowing with jewels and gold coins, with a big smile on his face. The image should be bright and colorful to capture the excitement and adventure of playing Gate of The Pharaohs."
$newMeta = "Our review of Gate of The Pharaohs covers gameplay, theme, and bonus features. Play for free on desktop or mobile devices."

$imagePromptPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Create a cartoon-style feature image*") {
        $imagePromptPara = $p
        break
    }
}

if ($imagePromptPara -ne $null) {
    $priorPara = $imagePromptPara.Previous()
    $insertionPoint = $d.Range($priorPara.Range.End, $priorPara.Range.End)

    $ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
    $xml = "<w:p $ns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Gate of The Pharaohs Slot for Free</w:t></w:r></w:p>" + `
           "<w:p $ns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$newMeta</w:t></w:r></w:p>"
    $insertionPoint.InsertXML($xml)

    # The InsertXML call above spliced the new italic run in immediately
    # before the original (now continuation) paragraph's own runs, so the
    # paragraph currently reads "<newMeta><oldPrompt>". Collapse that back
    # down to just the new text.
    $combined = $newMeta + $oldPrompt
    $d.Content.Find.Execute($combined, $false, $false, $false, $false, $false, $true, 1, $false, $newMeta, 2) | Out-Null
}
